$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived values for columns E..T (Ligand-expressing cells through
# Edge total expression derived specificity), rows 2..11.
$arr = New-Object 'object[,]' 10,16
$arr[0,0] = 1
$arr[0,1] = 0.3333333333333333
$arr[0,2] = 0.3642143333333334
$arr[0,3] = 1.092643
$arr[0,4] = 0.4800482050304226
$arr[0,5] = 0.4800482050304224
$arr[0,6] = 3
$arr[0,7] = 1
$arr[0,8] = 2.535712666666667
$arr[0,9] = 7.607138
$arr[0,10] = 0.04494879354621957
$arr[0,11] = 0.05070282964779482
$arr[0,12] = 0.9235428984148889
$arr[0,13] = 8.311886085734
$arr[0,14] = 0.02157758766014575
$arr[0,15] = 0.02433980236238719
$arr[1,0] = 1
$arr[1,1] = 0.3333333333333333
$arr[1,2] = 0.3642143333333334
$arr[1,3] = 1.092643
$arr[1,4] = 0.4800482050304226
$arr[1,5] = 0.4800482050304224
$arr[1,6] = 3
$arr[1,7] = 1
$arr[1,8] = 32.24261766666667
$arr[1,9] = 96.72785300000001
$arr[1,10] = 0.5715421877013505
$arr[1,11] = 0.6447070965264385
$arr[1,12] = 11.74322349838656
$arr[1,13] = 105.689011485479
$arr[1,14] = 0.2743678013051941
$arr[1,15] = 0.3094904844578921
$arr[2,0] = 1
$arr[2,1] = 0.3333333333333333
$arr[2,2] = 0.3642143333333334
$arr[2,3] = 1.092643
$arr[2,4] = 0.4800482050304226
$arr[2,5] = 0.4800482050304224
$arr[2,6] = 3
$arr[2,7] = 1
$arr[2,8] = 1.538811333333333
$arr[2,9] = 4.616434
$arr[2,10] = 0.02727742533206951
$arr[2,11] = 0.03076929413956839
$arr[2,12] = 0.5604571438957778
$arr[2,13] = 5.044114295062
$arr[2,14] = 0.01309447906851135
$arr[2,15] = 0.0147707444217529
$arr[3,0] = 1
$arr[3,1] = 0.3333333333333333
$arr[3,2] = 0.3642143333333334
$arr[3,3] = 1.092643
$arr[3,4] = 0.4800482050304226
$arr[3,5] = 0.4800482050304224
$arr[3,6] = 2
$arr[3,7] = 1
$arr[3,8] = 19.206297
$arr[3,9] = 38.412594
$arr[3,10] = 0.3404565075487166
$arr[3,11] = 0.2560262755732715
$arr[3,12] = 6.995208657657001
$arr[3,13] = 41.971251945942
$arr[3,14] = 0.1634355353396879
$arr[3,15] = 0.1229049540295733
$arr[4,0] = 1
$arr[4,1] = 0.3333333333333333
$arr[4,2] = 0.3642143333333334
$arr[4,3] = 1.092643
$arr[4,4] = 0.4800482050304226
$arr[4,5] = 0.4800482050304224
$arr[4,6] = 3
$arr[4,7] = 1
$arr[4,8] = 0.8899256666666666
$arr[4,9] = 2.669777
$arr[4,10] = 0.0157750858716439
$arr[4,11] = 0.01779450411292666
$arr[4,12] = 0.3241236834012222
$arr[4,13] = 2.917113150611
$arr[4,14] = 0.007572801656883434
$arr[4,15] = 0.008542219758816912
$arr[5,0] = 1
$arr[5,1] = 0.3333333333333333
$arr[5,2] = 0.3944893333333333
$arr[5,3] = 1.183468
$arr[5,4] = 0.5199517949695774
$arr[5,5] = 0.5199517949695774
$arr[5,6] = 3
$arr[5,7] = 1
$arr[5,8] = 2.535712666666667
$arr[5,9] = 7.607138
$arr[5,10] = 0.04494879354621957
$arr[5,11] = 0.05070282964779482
$arr[5,12] = 1.000311599398222
$arr[5,13] = 9.002804394584
$arr[5,14] = 0.02337120588607382
$arr[5,15] = 0.02636302728540762
$arr[6,0] = 1
$arr[6,1] = 0.3333333333333333
$arr[6,2] = 0.3944893333333333
$arr[6,3] = 1.183468
$arr[6,4] = 0.5199517949695774
$arr[6,5] = 0.5199517949695774
$arr[6,6] = 3
$arr[6,7] = 1
$arr[6,8] = 32.24261766666667
$arr[6,9] = 96.72785300000001
$arr[6,10] = 0.5715421877013505
$arr[6,11] = 0.6447070965264385
$arr[6,12] = 12.71936874824489
$arr[6,13] = 114.474318734204
$arr[6,14] = 0.2971743863961563
$arr[6,15] = 0.3352166120685464
$arr[7,0] = 1
$arr[7,1] = 0.3333333333333333
$arr[7,2] = 0.3944893333333333
$arr[7,3] = 1.183468
$arr[7,4] = 0.5199517949695774
$arr[7,5] = 0.5199517949695774
$arr[7,6] = 3
$arr[7,7] = 1
$arr[7,8] = 1.538811333333333
$arr[7,9] = 4.616434
$arr[7,10] = 0.02727742533206951
$arr[7,11] = 0.03076929413956839
$arr[7,12] = 0.6070446570124444
$arr[7,13] = 5.463401913112
$arr[7,14] = 0.01418294626355816
$arr[7,15] = 0.01599854971781548
$arr[8,0] = 1
$arr[8,1] = 0.3333333333333333
$arr[8,2] = 0.3944893333333333
$arr[8,3] = 1.183468
$arr[8,4] = 0.5199517949695774
$arr[8,5] = 0.5199517949695774
$arr[8,6] = 2
$arr[8,7] = 1
$arr[8,8] = 19.206297
$arr[8,9] = 38.412594
$arr[8,10] = 0.3404565075487166
$arr[8,11] = 0.2560262755732715
$arr[8,12] = 7.576679299332
$arr[8,13] = 45.460075795992
$arr[8,14] = 0.1770209722090287
$arr[8,15] = 0.1331213215436982
$arr[9,0] = 1
$arr[9,1] = 0.3333333333333333
$arr[9,2] = 0.3944893333333333
$arr[9,3] = 1.183468
$arr[9,4] = 0.5199517949695774
$arr[9,5] = 0.5199517949695774
$arr[9,6] = 3
$arr[9,7] = 1
$arr[9,8] = 0.8899256666666666
$arr[9,9] = 2.669777
$arr[9,10] = 0.0157750858716439
$arr[9,11] = 0.01779450411292666
$arr[9,12] = 0.3510661829595555
$arr[9,13] = 3.159595646636
$arr[9,14] = 0.008202284214760467
$arr[9,15] = 0.009252284354109745

$ws.Range("E2:T11").Value = $arr
